$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect temporarily to edit cell values, then restore
$ws.Unprotect()

# Update the confidential disclaimer text (date change 2021-07-13 -> 2021-07-14)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."
# Re-fit row 11 so the multi-line text doesn't leave a stray custom row height
$ws.Rows.Item(11).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5307654371401141
$ws.Range("E2").Value = -0.0003376553214479161

$ws.Range("D3").Value = 0.2695825604848546
$ws.Range("E3").Value = 0.00320769847634339

$ws.Range("D4").Value = 0.04979008296845878
$ws.Range("E4").Value = -0.009701492537313561

$ws.Range("D5").Value = 0.09377019696095083
$ws.Range("E5").Value = -0.002305475504322696

$ws.Range("D6").Value = 0.02685046385921559
$ws.Range("E6").Value = -0.01135029354207417

$ws.Range("D7").Value = 0.02924125858640602
$ws.Range("E7").Value = -0.01270696958028494

$ws.Range("E8").Value = -0.0006900276461444532

# Re-apply sheet protection to restore original protected state
$ws.Protect()
